# Add the new "2022-Q4" sheet (quarterly fund-holdings detail), insert it right
# after "总计" and before the existing "2022-Q3" sheet, and update the "总计"
# (summary) sheet with the new quarter's totals.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new sheet, positioned before "2022-Q3", and rename it.
# ------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Reference sheet used purely to clone cell formatting (bold header / bordered
# index column) that already exists elsewhere in the workbook.
$fmtSrc = $wb.Worksheets.Item("2022-Q3")

# ------------------------------------------------------------------
# 2. Populate the header row.
# ------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $q4.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}
$fmtSrc.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Populate the fund rows (code, name, size, position, ratio, value, rank).
#    "$null" marks the one cell (G26) that is a genuine numeric 0 rather than
#    a text value in the source data.
# ------------------------------------------------------------------
$q4rows = @(
  @("506003", "富国科创板两年定期开放混合", "13.76", "98.91", "6.62", "0.9109", 1),
  @("519026", "海富通中小盘混合", "16.23", "92.68", "3.57", "0.5794", 10),
  @("519029", "华夏稳增混合", "9.01", "94.55", "4.17", "0.3757", 10),
  @("001542", "国泰互联网+股票", "8.50", "92.88", "3.29", "0.2796", 9),
  @("007345", "富国科技创新灵活配置混合", "8.84", "94.61", "2.68", "0.2369", 9),
  @("000742", "国泰新经济灵活配置混合A", "5.03", "86.40", "3.30", "0.1660", 7),
  @("008009", "华商高端装备制造股票A", "6.14", "92.01", "2.64", "0.1621", 10),
  @("005819", "国泰优势行业混合A", "2.24", "90.14", "3.52", "0.0788", 6),
  @("506009", "国泰科创板两年定期开放混合", "2.17", "87.42", "3.30", "0.0716", 5),
  @("010912", "国泰成长价值混合A", "2.13", "86.46", "3.31", "0.0705", 7),
  @("000609", "华商新量化灵活配置混合A", "2.87", "85.51", "2.42", "0.0695", 10),
  @("012411", "海富通成长领航混合C", "1.71", "92.50", "3.92", "0.0670", 10),
  @("012410", "海富通成长领航混合A", "1.65", "92.50", "3.92", "0.0647", 10),
  @("010642", "农银汇理瑞祥一年持有期混合", "2.48", "21.31", "1.77", "0.0439", 6),
  @("001723", "华商新动力灵活配置混合", "0.75", "84.72", "4.96", "0.0372", 1),
  @("000530", "招商丰盛稳定增长灵活配置混合A", "0.55", "94.98", "4.53", "0.0249", 5),
  @("015585", "国泰优势行业混合C", "0.70", "90.14", "3.52", "0.0246", 6),
  @("000166", "中海信息产业精选混合", "0.72", "83.53", "3.25", "0.0234", 9),
  @("002417", "招商丰盛稳定增长灵活配置混合C", "0.17", "94.98", "4.53", "0.0077", 5),
  @("005997", "天弘裕利灵活配置混合C", "0.48", "46.59", "1.14", "0.0055", 5),
  @("014989", "国泰新经济灵活配置混合C", "0.14", "86.40", "3.30", "0.0046", 7),
  @("010913", "国泰成长价值混合C", "0.11", "86.46", "3.31", "0.0036", 7),
  @("002388", "天弘裕利灵活配置混合A", "0.09", "46.59", "1.14", "0.0010", 5),
  @("016050", "华商高端装备制造股票C", "0.01", "92.01", "2.64", "0.0003", 10),
  @("016048", "华商新量化灵活配置混合C", "0.00", "85.51", "2.42", $null, 10)
)

$r = 2
foreach ($row in $q4rows) {
    # Column A: plain 0-based row index, numeric.
    $q4.Cells.Item($r, 1).Value = ($r - 2)

    # Column B/C: fund code & name - always text.
    $q4.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]

    # Columns D/E/F: decimal-looking figures stored as text in the source.
    $q4.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4.Cells.Item($r, 6).Value = "'" + $row[4]

    # Column G: text, except the single genuine-zero numeric case.
    if ($row[5] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = "'" + $row[5]
    }

    # Column H: numeric rank.
    $q4.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# Clear the quote-prefix styling picked up from the "'" text markers above and
# apply the bordered/bold look used for the index column everywhere else.
$q4.Range("B2:G26").Style = "Normal"
$fmtSrc.Range("A2").Copy()
$q4.Range("A2:A26").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4. Update the "总计" summary sheet: insert the new "2022-Q4" row above the
#    existing quarters and renumber the index column.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# The freshly inserted row inherits the header's bold formatting; clear it
# back to plain/normal before filling in values (matches rows below it).
$total.Range("B2:D2").Style = "Normal"

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 25
$total.Cells.Item(2, 4).Value = 3.31

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

# Match the index-column formatting used by the rest of the sheet.
$total.Cells.Item(3, 1).Copy()
$total.Range("A2").PasteSpecial(-4122)

# Restore the originally-active tab ("2022-Q1") - adding/renaming sheets above
# shifts the active selection, and this wasn't part of the intended change.
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Output "ok"
